$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 841; this shifts the existing rows 841:882
# down to 842:883 and grows the sheet's used range to row 883.
$ws.Rows.Item(841).Insert()

# Fill the new row 841 with the inserted data point (2026/02/21, 土, 7, 201).
# Column A/B are copied from row 840 (which already holds the same date
# "2026/02/21" / weekday "土" text) instead of assigning literal strings,
# so Excel's text-to-date autoconversion never kicks in and no new
# number-format style gets minted in the process.
$ws.Range("A840:B840").Copy($ws.Range("A841:B841"))
$ws.Range("C841").Value = 7
$ws.Range("D841").Value = 201
